# Carparks_Keep_Throw_List.xlsx
# "[Q3] Update list to contain only carparks to keep"
#
# The sheet has a header row (Carpark | Keep/Throw) followed by one row per
# carpark, each flagged "Keep" or "Throw" in column B. This edit removes every
# "Throw" row so only the carparks that are being kept remain, then turns on
# an AutoFilter for the (now constant) Keep/Throw column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row, then scan bottom-up deleting any row whose
# Keep/Throw column (B) reads "Throw". Scanning bottom-up means row indices
# of not-yet-visited rows are never disturbed by the deletions above them.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = $lastRow; $r -ge 2; $r--) {
    $flag = $ws.Cells.Item($r, 2).Value()
    if ($flag -eq "Throw") {
        $ws.Rows.Item($r).Delete()
    }
}

# Recompute how many rows are left (header + remaining "Keep" carparks) and
# apply an AutoFilter over the Keep/Throw column for that range.
$lastRow = $ws.UsedRange.Rows.Count
$filterRef = "B1:B" + $lastRow
$ws.Range($filterRef).AutoFilter()

# Excel tracks the active AutoFilter range as a hidden workbook-scoped
# defined name.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$B`$1:`$B`$" + $lastRow)
$filterName.Visible = $false

# Leave the view scrolled near the bottom of the (now shorter) list, with
# E29 as the active cell, matching the saved view state.
$ws.Range("E29").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 144
$win.ScrollColumn = 1
